# ---------------------------------------------------------------------------
# Edit: Update a driver installation manual
#
# 1) In the "=====" separator paragraph, drop its paragraph-mark formatting
#    (the <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>) and
#    remove the _GoBack bookmark that lived there.
# 2) Replace the final empty paragraph (right before the sectPr) with a page
#    break followed by a whole new "APPENDIX A." section, re-homing the
#    _GoBack bookmark onto the new appendix heading.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------------
# Remove the _GoBack bookmark wherever it currently lives.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Find the "====...====" separator paragraph and strip the paragraph-mark's
# run formatting (the eastAsia rFonts hint living in w:pPr/w:rPr) by merging
# the paragraph mark into the next paragraph and then re-splitting: the new
# paragraph mark created by InsertParagraphAfter carries no formatting.
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text -like "*=================*") {
    $endOfPara = $p.Range.End
    $markRange = $d.Range($endOfPara - 1, $endOfPara)
    $markRange.Delete()
    $splitPoint = $d.Range($endOfPara - 1, $endOfPara - 1)
    $splitPoint.InsertParagraphAfter()
    break
  }
}

# --- Edit 2 -------------------------------------------------------------------
# Replace the final (empty) paragraph before the sectPr with a page break and
# the new Appendix A content.
$last = $d.Paragraphs.Last
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:widowControl/>
    <w:jc w:val="left"/>
  </w:pPr>
  <w:r>
    <w:br w:type="page"/>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:b/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>APPENDIX A.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>setup for handling COPPER and TTRX drivers</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:tab/>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Sep. 2, 2014 S. Yamada</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>The following is a procedure to set up</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> COPPER and TTRX drivers in ecl02</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>1, copy *.ko and script files from ecl01</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>2, make a direcotry for module files</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>[b2daq@ecl02:yamadas]$ sudo mkdir /tftpboot/copper/root/lib/modules/2.6.18/misc</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>3, copy module files in the misc directory</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>[b2daq@ecl02:yamadas]$ ls *.ko</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>copper.0.1.5.ko              copper.131010.ko           cprfin_fngeneric.131008.ko  vme_universe.ko</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>copper.0.1.6.chksumcheck.ko  copper.ko                  cprfin_fngeneric.ko</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>copper.0.1.6.ko              cprfin_fngeneric.0.1.5.ko  ttrx_fifo.ko</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>copper.0.1.6.timeout1s.ko    cprfin_fngeneric.0.1.6.ko  ttrx.ko</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>[b2daq@ecl02:yamadas]$ sudo cp *.ko /tftpboot/copper/root/lib/modules/2.6.18/misc</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>4, copy scripts</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>[b2daq@ecl02:yamadas]$ sudo cp copper /tftpboot/copper/root/etc/rc.d/init.d/</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>[b2daq@ecl02:yamadas]$ sudo cp ttrx /tftpboot/copper/root/etc/rc.d/init.d/</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>5, Add copper and ttrx in the startup list</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> in COPPER</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>[root@cpr5014:b2daq]#  sudo /sbin/chkconfig --add copper</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve">[root@cpr5014:b2daq]#  sudo /sbin/chkconfig --add ttrx  </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>6, reboot COPPERs or install the drivers by hand</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>[root@cpr5014:b2daq]# /sbin/service copper start</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>Loading COPPER driver:                                     [  OK  ]</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>Loading FINESSE driver:                                    [  OK  ]</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve">[root@cpr5014:b2daq]# /sbin/service ttrx start  </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>Loading TT-RX device driver:                               [  OK  ]</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>Loading TT-RX FIFO driver:</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$last.Range.InsertXML($xml)
